# Update VIC second-doses daily series to 21 December:
# insert 4 new leading rows (17-20 Dec data) above the existing series,
# pushing all prior rows down by 4, then set the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows above the first data row, shifting everything down.
$ws.Rows("2:5").Insert()

# Copy the (now-shifted) first original data row's formatting (date style /
# number style + row height) onto the newly inserted rows so the new rows
# match the rest of the table.
$ws.Range("A6:B6").Copy()
$ws.Range("A2:B5").PasteSpecial(-4122)
$ws.Range("A2:B5").RowHeight = 18

# New daily figures (date serial, cumulative second doses), newest first.
$ws.Range("A2").Value = 44550
$ws.Range("B2").Value = 5237300
$ws.Range("A3").Value = 44549
$ws.Range("B3").Value = 5236438
$ws.Range("A4").Value = 44548
$ws.Range("B4").Value = 5234686
$ws.Range("A5").Value = 44547
$ws.Range("B5").Value = 5230468

# Match the saved selection from the source commit.
$ws.Range("G8").Select()
